# Add the latest Tesla stock-price sample to the tracker sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$row = 23
$rng = $ws.Range("A" + $row + ":C" + $row)

# Force text formatting on the text columns first so Excel does not
# auto-convert the date / currency-looking strings into a date serial
# number or a formatted currency number.
$rng.NumberFormat = "@"

$ws.Cells.Item($row, 1).Value = "2025-05-12"
$ws.Cells.Item($row, 2).Value = "11:23:03"
$ws.Cells.Item($row, 3).Value = "$298.26"
$ws.Cells.Item($row, 4).Value = 298.26

# Drop the temporary "Text" number format again so the new cells end up
# unstyled (matching the rest of the data rows), while keeping the
# literal string values that were just entered.
$rng.ClearFormats()
